$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# Fill in category (C) and price (D) for rows 43-45, which were
# previously blank placeholder rows in the table.
$ws.Range("C43").Value = "Munich"
$ws.Range("D43").Value = 25
$ws.Range("C44").Value = "Kioto"
$ws.Range("D44").Value = 30
$ws.Range("C45").Value = "Huesca"
$ws.Range("D45").Value = 25

# Rows 46-48 were placeholder rows with only a "tipo" value set; they
# become fully empty (keeping their row height only).
$ws.Range("A46:I48").Clear()

# The table now only spans down to row 45.
$lo = $ws.ListObjects.Item("Table_1")
$lo.Resize($ws.Range("A1:I45"))

# Data validation ranges shrink along with the table.
$ws.Range("C2:C47").Validation.Delete()
$dvCategoria = $ws.Range("C2:C45").Validation
$dvCategoria.Add(3, 1, 1, "=datos!`$C`$2:`$C`$14")
$dvCategoria.ShowInput = $false
$dvCategoria.ShowError = $true
$dvCategoria.IgnoreBlank = $true

$ws.Range("B2:B47").Validation.Delete()
$dvTipo = $ws.Range("B2:B45").Validation
$dvTipo.Add(3, 1, 1, "=datos!`$A`$2:`$A`$5")
$dvTipo.ShowInput = $false
$dvTipo.ShowError = $true
$dvTipo.IgnoreBlank = $true

# Trim the trailing empty filler rows that are no longer needed.
$ws.Rows("1004:1006").Delete()
